# Update "想去人数" (interest-count) figures in column F across all four
# sheets, matching the freshly re-scraped data snapshot (gh-pages output
# generated at 456a3b4). Only column F numeric values change; everything
# else in the workbook stays the same.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1217
$ws1.Range("F4").Value = 1279
$ws1.Range("F6").Value = 179
$ws1.Range("F7").Value = 556
$ws1.Range("F8").Value = 20
$ws1.Range("F9").Value = 347
$ws1.Range("F11").Value = 1271
$ws1.Range("F12").Value = 29197
$ws1.Range("F13").Value = 4336
$ws1.Range("F14").Value = 43
$ws1.Range("F15").Value = 266
$ws1.Range("F16").Value = 492
$ws1.Range("F17").Value = 47
$ws1.Range("F19").Value = 15
$ws1.Range("F21").Value = 347
$ws1.Range("F22").Value = 636
$ws1.Range("F23").Value = 279
$ws1.Range("F24").Value = 285
$ws1.Range("F25").Value = 359
$ws1.Range("F27").Value = 73
$ws1.Range("F29").Value = 670
$ws1.Range("F30").Value = 217
$ws1.Range("F31").Value = 104
$ws1.Range("F32").Value = 554
$ws1.Range("F33").Value = 80
$ws1.Range("F35").Value = 646
$ws1.Range("F36").Value = 245
$ws1.Range("F37").Value = 41

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 898
$ws2.Range("F8").Value = 1
$ws2.Range("F11").Value = 276
$ws2.Range("F12").Value = 4259
$ws2.Range("F23").Value = 4252

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 300
$ws3.Range("F3").Value = 265
$ws3.Range("F4").Value = 1209

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 898
$ws4.Range("F10").Value = 1217
$ws4.Range("F11").Value = 1279
$ws4.Range("F12").Value = 179
$ws4.Range("F13").Value = 556
$ws4.Range("F14").Value = 20
$ws4.Range("F15").Value = 347
$ws4.Range("F18").Value = 1271
$ws4.Range("F21").Value = 276
$ws4.Range("F28").Value = 492
$ws4.Range("F29").Value = 47
$ws4.Range("F30").Value = 15
$ws4.Range("F34").Value = 347
$ws4.Range("F35").Value = 636
$ws4.Range("F36").Value = 279
$ws4.Range("F38").Value = 73
$ws4.Range("F40").Value = 670
$ws4.Range("F42").Value = 217
$ws4.Range("F43").Value = 104
$ws4.Range("F46").Value = 80
$ws4.Range("F48").Value = 646
$ws4.Range("F49").Value = 245
$ws4.Range("F50").Value = 41
